$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column AN (rows 1-11) to column AO so the new
# column inherits the same styles (header style + centered integer style).
$ws.Range("AN1:AN11").Copy()
$ws.Range("AO1:AO11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AO1").Value = "03-ago"
$ws.Cells.Item(2, 41).Value = 16
$ws.Cells.Item(3, 41).Value = 14
$ws.Cells.Item(4, 41).Value = 11
$ws.Cells.Item(5, 41).Value = 17
$ws.Cells.Item(6, 41).Value = 8
$ws.Cells.Item(7, 41).Value = 18
$ws.Cells.Item(8, 41).Value = 16
$ws.Cells.Item(9, 41).Value = 16
$ws.Cells.Item(10, 41).Value = 13
$ws.Cells.Item(11, 41).Value = 12

$ws.Range("AT10").Select() | Out-Null
